# Governance Body Suite.xlsx — 20-Feb-2017 client-data update.
#
# The "GBCreation" sheet drives a test fixture row (row 3) used to seed a
# governance-body creation test case. The client reference in that row is
# being swapped from "ABC News" to "HSBC":
#   A3: supplier            "ABC News"                             -> "HSBC"
#   B3: contract             "Master Service Agreement - ABC News" -> "Master Service Agreement - HSBC"
#
# (All other cells/strings in the workbook keep their values; any shared-
# string index renumbering elsewhere is just a side effect of the string
# table no longer needing the old "ABC News" / "Master Service Agreement -
# ABC News" entries.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GBCreation")

$ws.Range("A3").Value = "HSBC"
$ws.Range("B3").Value = "Master Service Agreement - HSBC"
